$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''30.501.43'
$ws.Range('E2').Value = '  -1.01%  '
$ws.Range('D3').Value = '''2.111.21'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '''334.72'
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = '''0.5247'
$ws.Range('E7').Value = '  -1.83%  '
$ws.Range('D8').Value = '''0.4529'
$ws.Range('E8').Value = '  +2.81%  '
$ws.Range('D9').Value = '''53.60'
$ws.Range('E9').Value = '  +13.73%  '
$ws.Range('D10').Value = '''0.09006'
$ws.Range('E10').Value = '  +0.05%  '
$ws.Range('E11').Value = '  -1.43%  '
$ws.Range('E12').Value = '  -1.89%  '
$ws.Range('D13').Value = '''2.106.44'
$ws.Range('E13').Value = '  -0.56%  '
$ws.Range('D14').Value = '''6.791'
$ws.Range('E14').Value = '  +0.35%  '
$ws.Range('D15').Value = '''7.825'
$ws.Range('E15').Value = '  +0.15%  '
$ws.Range('D16').Value = '''96.83'
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('D17').Value = '''1.003'
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').Value = '''0.00001128'
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('D19').Value = '''0.06625'
$ws.Range('E19').Value = '  -0.70%  '
$ws.Range('D20').Value = '''19.34'
$ws.Range('E20').Value = '  +0.92%  '
$ws.Range('D22').Value = '''6.316'
$ws.Range('E22').Value = '  -0.49%  '
$ws.Range('D23').Value = '''30.557.85'
$ws.Range('E23').Value = '  -1.03%  '
$ws.Range('D24').Value = '''12.42'
$ws.Range('E24').Value = '  +0.62%  '
$ws.Range('D25').Value = '''2.344'
$ws.Range('E25').Value = '  +2.38%  '
$ws.Range('D26').Value = '''2.349.90'
$ws.Range('E26').Value = '  -0.71%  '
$ws.Range('D27').Value = '''22.41'
$ws.Range('E27').Value = '  -1.59%  '
$ws.Range('D28').Value = '''2.586'
$ws.Range('E28').Value = '  -0.41%  '
$ws.Range('D29').Value = '''163.50'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').Value = '''132.77'
$ws.Range('E30').Value = '  -0.40%  '
$ws.Range('D31').Value = '''1.204'
$ws.Range('E31').Value = '  +2.09%  '
$ws.Range('E32').Value = '  -0.75%  '
$ws.Range('D33').Value = '''1.663'
$ws.Range('E33').Value = '  +7.52%  '
$ws.Range('D34').Value = '''6.169'
$ws.Range('E34').Value = '  -1.26%  '
$ws.Range('D35').Value = '''3.940'
$ws.Range('E35').Value = '  -1.94%  '
$ws.Range('D36').Value = '''10.66'
$ws.Range('E36').Value = '  +12.34%  '
$ws.Range('D37').Value = '''0.02579'
$ws.Range('E37').Value = '  -0.98%  '
$ws.Range('D38').Value = '''0.06836'
$ws.Range('E38').Value = '  +1.05%  '
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('D40').Value = '''12.81'
$ws.Range('E40').Value = '  -0.53%  '
$ws.Range('D41').Value = '''0.2298'
$ws.Range('E41').Value = '  +0.38%  '
$ws.Range('D42').Value = '''0.6942'
$ws.Range('E42').Value = '  +1.15%  '
$ws.Range('D43').Value = '''1.245'
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('D44').Value = '''2.407'
$ws.Range('E44').Value = '  +8.05%  '
$ws.Range('E45').Value = '  +0.24%  '
$ws.Range('D46').Value = '''0.6418'
$ws.Range('E46').Value = '  -0.64%  '
$ws.Range('D47').Value = '''14.07'
$ws.Range('E47').Value = '  -0.73%  '
$ws.Range('D48').Value = '''3.661'
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('D49').Value = '''1.250'
$ws.Range('E49').Value = '  -2.16%  '
$ws.Range('D50').Value = '''1.215'
$ws.Range('E50').Value = '  +3.85%  '
$ws.Range('D51').Value = '''83.35'
$ws.Range('E51').Value = '  +0.22%  '
